$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran the simulation for the first convergence block (rows 4-8): only the
# "diagonal" raw measurements (C4, D5, E6, F7, G8) are populated with the new
# results; the rest of that block is cleared out. Downstream ratio/LOG
# formulas recalc automatically (producing #DIV/0!/#NUM! where an operand is
# now blank), matching the "Fixed error in F" commit.

$ws.Range("C4").Value = 0.00005475780041680050
$ws.Range("D4:G4").ClearContents()

$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 0.00004410779576484220
$ws.Range("E5:G5").ClearContents()

$ws.Range("C6:D6").ClearContents()
$ws.Range("E6").Value = 0.00004009327123032160
$ws.Range("F6:G6").ClearContents()

$ws.Range("C7:E7").ClearContents()
$ws.Range("F7").Value = 0.00003869669795327120
$ws.Range("G7").ClearContents()

$ws.Range("C8:F8").ClearContents()
$ws.Range("G8").Value = 0.00003822580448587120

# Restore the active cell/selection left by the author after the edit.
$ws.Range("D25").Select()
